$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.695.20"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.606.48"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.06"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.91%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.837.16"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "1.633.18"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "29.720.26"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.11"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.40"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").Value = "0.0₃0698"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.07"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.46"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").Value = "1.427.75"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.549"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.92"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.25%  "
$ws.Range("E42").Value = "  +5.81%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.95"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.817"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.980"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +16.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.36"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "1.746.05"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.69"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "0.0₆0106"
$ws.Range("E51").Value = "  +5.57%  "
